$p = $ppt.ActivePresentation
try {
    $r = $ppt.OpenThemeFile("foo.thmx", "notesMaster")
    Write-Output ("result=" + $r)
} catch {
    Write-Output ("ERR: " + $_.Exception.Message)
}
